# Codebook index sheet maintenance edit:
# - removes three now-unused spacer rows (old rows 24, 25, 28) from the
#   "index" sheet, shifting the D5_IR / D5_results / D5_Table_3_IR /
#   D5_Figure_IR blocks up
# - adds two new "source" references in column F for the D5_results and
#   D5_Figure_IR rows
# - updates the frozen-pane scroll position / current selection
# - shrinks the conditional-formatting "no duplicates" range to start at
#   the new row 26 instead of the old row 29

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the three blank spacer rows -----------------------------------
# Old layout (rows 24-30):
#   24: (blank C/D only)
#   25: (blank B only)
#   26: D5_IR / incidence rate... / 08_T4_30_create_D5_IR
#   27: D5_results_from_analysis / 08_T4_40_analysis
#   28: (blank B only)
#   29: D5_Table_3_IR / D4_Cube_persontime_bleeding D3_source_population
#   30: D5_Figure_IR / 09_T5_35_create_Figure_1
# After removing 24, 25, 28 the remaining rows shift up to become the new
# rows 24-27.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(26).Delete()

# --- Fill in the new column-F references on the (now) row 25 and 27 -------
$ws.Range("F25").Value = "D4_analytic_dataset"
$ws.Range("F27").Value = "D5_IR"

# --- Shrink / move the "no duplicates" conditional formatting -------------
# It used to start at B29 (and excluded C21:C23 while covering a separate
# C25:C28 block); now that those rows are gone it simply starts at B26 and
# covers a single contiguous C21:C25 block instead.
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -ne '$C$19') {
        $newRange = $ws.Range("B26:C1048576,B1:C4,B12:C18,B7:C10,B5:B6,B11,C21:C25")
        $fc.ModifyAppliesToRange($newRange)
        $fc.Priority = 5
    }
}

# --- Update the frozen-pane view / selection -------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 5
$ws.Rows.Item(26).Select()
